# Fix computer science only filter: restore publications that a prior
# (too-aggressive) CS-only filter had dropped, re-inserting them in their
# correct chronological position and renumbering the trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has 23 data rows (A1:F23); it needs to grow to 27
# (A1:F27), i.e. 4 new rows. Insert them ahead of the final row (current
# row 23, "DanXe...") - row-by-row copies of row 22 seed each new row with
# matching formatting (border/bold index-column style, empty placeholder
# cell in column E) which we then overwrite with the real values.
$ws.Range("A23:F26").EntireRow.Insert()
$ws.Range("A22:F22").Copy($ws.Range("A23:F23"))
$ws.Range("A22:F22").Copy($ws.Range("A24:F24"))
$ws.Range("A22:F22").Copy($ws.Range("A25:F25"))
$ws.Range("A22:F22").Copy($ws.Range("A26:F26"))

# Column F ("Year") stores its values as text in this workbook (not
# numbers) - force text format before writing so the new/changed cells
# keep the same storage type as every other row.
$ws.Range("F20:F27").NumberFormat = "@"

# Column E ("Quartile") is blank text (present, empty-string cell, not a
# truly-blank cell) for every non-journal-article row. A leading
# apostrophe forces the engine to store it the same way (empty Text)
# instead of collapsing an assigned "" into a blank cell.
$ws.Range("E23").Value = "'"
$ws.Range("E24").Value = "'"
$ws.Range("E25").Value = "'"
$ws.Range("E27").Value = "'"

# Row 20: was "Flying in XR...", now becomes "Fashion in the Metaverse..."
$ws.Range("B20").Value = "Fashion in the Metaverse: Technologies, Applications, and Opportunities"
$ws.Range("D20").Value = "International Conference on Entertainment Computing"
$ws.Range("F20").Value = "2023"

# Row 21: was "WiXaRd...", now becomes "Social cognition and Metaverse..."
$ws.Range("B21").Value = "Social cognition and Metaverse: understanding and disposition towards Human Digital Twins"
$ws.Range("D21").Value = "Atti del Congresso dell'Associazione Italiana di Psicologia 2022"
$ws.Range("F21").Value = "2022"

# Row 22: was "M-AGEW...", now becomes "Flying in XR..." (restored)
$ws.Range("B22").Value = "Flying in XR: Bridging Desktop Applications in eXtended Reality through Deep Learning"
$ws.Range("D22").Value = "IEEE Conference on Virtual Reality and 3D User Interfaces Abstracts and Workshops (VRW) (pp. 264-272)"
$ws.Range("F22").Value = "2024"

# Row 23 (new): "WiXaRd..." (restored)
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "WiXaRd: Towards a holistic distributed platform for multi-party and cross-reality WebXR experiences"
$ws.Range("C23").Value = "inproceedings"
$ws.Range("D23").Value = "IEEE Conference on Virtual Reality and 3D User Interfaces Abstracts and Workshops (VRW)."
$ws.Range("F23").Value = "2024"

# Row 24 (new): "M-AGEW..." (restored)
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "M-AGEW: Empowering Outdoor Workouts with Data-Driven Augmented Reality Assistance"
$ws.Range("C24").Value = "inproceedings"
$ws.Range("D24").Value = "2024 IEEE International Conference on Artificial Intelligence and eXtended and Virtual Reality (AIxVR)"
$ws.Range("F24").Value = "2024"

# Row 25 (new): "AI for Enhancing and Preserving Dance Cultural Heritage..."
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = "AI for Enhancing and Preserving Dance Cultural Heritage: a Case Study on Rudolf Nureyev's Costumes"
$ws.Range("C25").Value = "inproceedings"
$ws.Range("D25").Value = "The First International Conference on Artificial Intelligence and Immersive Virtual Reality - AIVR 2024"
$ws.Range("F25").Value = "2024"

# Row 26: the old "DanXe..." row, shifted down and renumbered.
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = "DanXe: An extended artificial intelligence framework to analyze and promote dance heritage"
$ws.Range("C26").Value = "article"
$ws.Range("D26").Value = "Digital Applications in Archaeology and Cultural Heritage"
$ws.Range("E26").Value = "Q2"
$ws.Range("F26").Value = "2024"

# Row 27 (new): "An Ethical Framework for Trustworthy Neural Rendering..."
# (this is the shifted-down original last row, which had a "Q2" Quartile
# value in column E - it was reset to blank text above since this new
# entry has no quartile)
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = "An Ethical Framework for Trustworthy Neural Rendering applied in Cultural Heritage and Creative Industries"
$ws.Range("C27").Value = "inproceedings"
$ws.Range("D27").Value = "CVPR 2024, AI for 3D Generation Workshop"
$ws.Range("F27").Value = "2024"
